# FA20_TestData_ManageTransactionsAssetInvoices_21C.xlsx
# "Add files via upload" - re-pointing this copy of the template away from the
# old Selenium/Oracle-Fusion login URL towards the shared IBM-implementation
# login, i.e. the explicit URL/user/password that used to live in Z2:AB2
# (with a live hyperlink on Z2) are wiped out.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")
$ws.Activate()

# Drop the hyperlink object that lived on Z2 (https://edrx.fa.us2.oraclecloud.com/)
# before clearing the cell text, otherwise the hyperlink annotation survives an
# empty cell.
$ws.Hyperlinks.Delete()

# Clear the URL / username / password cells entirely (they become blank cells,
# not just empty strings) - this also drops the now-unused shared strings
# "https://edrx.fa.us2.oraclecloud.com", "IBM_IMPLEMENTATION_USER" and
# "Oracle1234" from the workbook.
$ws.Range("Z2:AB2").ClearContents()

# Reflect the on-screen selection left behind after clearing that block.
$ws.Range("Z2:AB2").Select()
